$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "pop" column (column C), shifting n_pharmacies and n_items left
# one column over (C4->C3, E->D etc.)
[void]$ws.Range("C:C").Delete()

# Move the active selection to match where the editor ended up (D8)
[void]$ws.Range("D8").Select()
